# BA_tasks.xlsx — "feat: added convolutional AE"
#
# The sheet keeps a running "done notes" log in column D. The entry that used
# to say "try Convolutional AE: https://blog.keras.io/..." is removed from its
# slot; every later entry in the log shifts up by one row to fill the gap, and
# a brand-new final entry "Conv. AE is bad" is appended. The (unrelated)
# Convolutional-AE link text itself is folded into the end of the long
# "Eigenfaces: 3-4 components, ..." note instead. Two cells (D8 and D16) also
# pick up a new yellow highlight style, row heights for rows 22/26 grow to fit
# their new wrapped text, and the now-unused trailing rows (30-37) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "done notes" log in column D --------------------------------

$ws.Range('D8').Value = 'CNN Ähnlichkeiten auf images der PDFs erkennen lassen'
$ws.Range('D8').Interior.Color = 65535
$ws.Range('D8').WrapText = $true

$ws.Range('D9').Value = 'Google''s universal Sentence Encoder: Alter shapes to fix problem (HOW?), InferSent has same difference problem'

$ws.Range('D10').Value = 'Add pipeline image for different models, which shows if stopwords were removed by me or model etc. to bachelor thesis'

$ws.Range('D11').Value = 'analysis/ evaluation ideas cf. Notability 23.08.2023'

$ws.Range('D12').Value = 'why are cluster sizes of PCA results imbalanced? Because they are too sparse'

$ws.Range('D13').Value = 'Cluster: bewerten lassen'

$ws.Range('D14').Value = 'Cluster für Bilder und Text nutzbar'

$ws.Range('D15').Value = 'Universal Sentence Encoder: nicht auf GPU, finde Artikel mit gleichem Problem online'

$ws.Range('D16').Value = 'Universal Sentence Encoder: try huggingface version'
$ws.Range('D16').Interior.Color = 65535
$ws.Range('D16').WrapText = $true

$ws.Range('D17').Value = 'Eigenfaces: find paper'

$ws.Range('D18').Value = 'Eigenfaces: hohe dim <-> sparse <-> Clustern schwer'

$ws.Range('D20').Value = 'Why are residual graph so similar for both latent dimensions?'

$ws.Range('D22').Value = 'BA: Pipeline Bild s. Christians Zeichnung, 62GB Daten -> offline verarbeiten -> DB -> auf kleinem System durchsuchbar'

$ws.Range('D23').Value = 'BA: AE latent space normal verteilt???? '

$ws.Range('D26').Value = 'BA: Kmeans Hypothese ist, dass Daten normalverteilt. Deshalb Kmeans cluster eher rund -> deshalb wäre anderer Algo, z.B. Var. Bayesian Mixture Model interessant'

$ws.Range('D27').Value = 'BA: OPTICS etc. in Absatz Funktion erklären + Referenz'

$ws.Range('D28').Value = 'Lizenzen: GPL kann man nicht so einfach nutzen'

$ws.Range('D29').Value = 'Conv. AE is bad'

# D24 and D25 no longer carry a note (their old text moved up into D14/D15) -
# clear them completely (not just the contents) so no stray empty cell is left.
$ws.Range('D24').Clear()
$ws.Range('D25').Clear()

# --- Update the text that now includes the "tried Convolutional AE" mention -
$ws.Range('D51a_placeholder').Value = $null
